$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-08-15 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-16 Saturday", 2)

# Update the division problems in the table, cell by cell (row, column)
# so that duplicate cell texts are handled unambiguously by position.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "44÷7=6, 2"
$t.Cell(1, 2).Range.Text = "24÷8=3, 0"
$t.Cell(1, 3).Range.Text = "88÷3=29, 1"
$t.Cell(1, 4).Range.Text = "12÷6=2, 0"
$t.Cell(1, 5).Range.Text = "56÷4=14, 0"

$t.Cell(5, 1).Range.Text = "26÷5=5, 1"
$t.Cell(5, 2).Range.Text = "56÷9=6, 2"
$t.Cell(5, 3).Range.Text = "30÷3=10, 0"
$t.Cell(5, 4).Range.Text = "96÷3=32, 0"
$t.Cell(5, 5).Range.Text = "70÷3=23, 1"

$t.Cell(9, 1).Range.Text = "74÷6=12, 2"
$t.Cell(9, 2).Range.Text = "73÷2=36, 1"
$t.Cell(9, 3).Range.Text = "75÷3=25, 0"
$t.Cell(9, 4).Range.Text = "12÷8=1, 4"
$t.Cell(9, 5).Range.Text = "35÷7=5, 0"

$t.Cell(13, 1).Range.Text = "58÷6=9, 4"
$t.Cell(13, 2).Range.Text = "73÷5=14, 3"
$t.Cell(13, 3).Range.Text = "10÷3=3, 1"
$t.Cell(13, 4).Range.Text = "69÷7=9, 6"
$t.Cell(13, 5).Range.Text = "90÷9=10, 0"

$t.Cell(17, 1).Range.Text = "62÷3=20, 2"
$t.Cell(17, 2).Range.Text = "63÷2=31, 1"
$t.Cell(17, 3).Range.Text = "61÷6=10, 1"
$t.Cell(17, 4).Range.Text = "24÷9=2, 6"
$t.Cell(17, 5).Range.Text = "40÷3=13, 1"
